# Updated symbol list on Fri Jan 13 09:34:56 UTC 2023 with GitHub Actions
# Apply new Price / Volume(1h) figures to the cryptos sheet, keeping cell
# contents as plain text (matching the sheet's existing inline-string cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.12"
$ws.Range("E2").Value = "'1.27%"
$ws.Range("D3").Value = "'29.63"
$ws.Range("E3").Value = "'4.15%"
$ws.Range("D4").Value = "'5.126"
$ws.Range("E4").Value = "'1.42%"
$ws.Range("D5").Value = "'0.06705"
$ws.Range("E5").Value = "'3.21%"
$ws.Range("D6").Value = "'7.331"
$ws.Range("E6").Value = "'1.57%"
$ws.Range("D7").Value = "'3.397"
$ws.Range("E7").Value = "'0.96%"
$ws.Range("D8").Value = "'1.364"
$ws.Range("E8").Value = "'-2.05%"
$ws.Range("D9").Value = "'0.9189"
$ws.Range("E9").Value = "'0.14%"
$ws.Range("D10").Value = "'0.1593"
$ws.Range("E10").Value = "'3.41%"
$ws.Range("E11").Value = "'1.64%"
$ws.Range("D12").Value = "'0.07731"
$ws.Range("E12").Value = "'1.51%"
$ws.Range("D13").Value = "'0.02929"
$ws.Range("E13").Value = "'4.99%"
$ws.Range("D14").Value = "'0.08990"
$ws.Range("E14").Value = "'0.27%"
$ws.Range("D15").Value = "'0.001586"
$ws.Range("E15").Value = "'0.22%"
$ws.Range("D16").Value = "'0.04482"
$ws.Range("E16").Value = "'0.94%"
$ws.Range("D17").Value = "'0.0006452"
$ws.Range("E17").Value = "'1.84%"
$ws.Range("D18").Value = "'0.006281"
$ws.Range("E18").Value = "'1.88%"
$ws.Range("D19").Value = "'3.445"
$ws.Range("E19").Value = "'-0.20%"
$ws.Range("D20").Value = "'2.227"
$ws.Range("E20").Value = "'-0.82%"
$ws.Range("E21").Value = "'1.05%"
$ws.Range("E22").Value = "'-2.92%"
$ws.Range("D23").Value = "'4.070"
$ws.Range("E23").Value = "'1.49%"
$ws.Range("E24").Value = "'2.39%"
$ws.Range("E25").Value = "'0.88%"
$ws.Range("D26").Value = "'0.004123"
$ws.Range("E26").Value = "'-7.61%"
$ws.Range("D27").Value = "'0.0001199"
$ws.Range("E27").Value = "'-0.08%"
$ws.Range("E28").Value = "'-0.14%"
$ws.Range("D40").Value = "'0.04273"
$ws.Range("E40").Value = "'3.85%"
$ws.Range("D41").Value = "'0.006723"
$ws.Range("E41").Value = "'0.77%"
$ws.Range("E42").Value = "'0.58%"
$ws.Range("D43").Value = "'0.002167"
$ws.Range("E43").Value = "'5.76%"
$ws.Range("D44").Value = "'0.01198"
$ws.Range("E44").Value = "'3.68%"
$ws.Range("D45").Value = "'0.00005704"
$ws.Range("E45").Value = "'5.74%"
$ws.Range("E47").Value = "'-29.47%"

# Excel auto-applies a number/percent style to text that *looks* numeric;
# explicitly restore the default "Normal" style so only the text content changes.
$touched = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","E21","E22","D23","E23","E24","E25","D26","E26","D27","E27","E28","D40","E40","D41","E41","E42","D43","E43","D44","E44","D45","E45","E47")
foreach ($addr in $touched) {
    $ws.Range($addr).Style = "Normal"
}
